# #11 updated menu Mock up
#
# The footer "Date" placeholder on the slide master and on every slide
# layout shows a fixed date field (type="datetimeFigureOut") that reads
# 11/25/2021. Bump it to 11/27/2021 everywhere it appears (slide master
# + all custom layouts), mirroring what PowerPoint does when the fixed
# date shown in Insert > Header & Footer is updated.

$p = $ppt.ActivePresentation
$newDate = "11/27/2021"

# ppPlaceholderDate = 16
$ppPlaceholderDate = 16

function Update-DateShapes($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $phType = $null
        try { $phType = $shp.PlaceholderFormat.Type } catch { $phType = $null }

        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -ne $text) {
                $shp.TextFrame.TextRange.Text = $text
            }
        }
    }
}

# 1) Slide master's Date placeholder
$master = $p.SlideMaster
Update-DateShapes $master.Shapes $newDate

# 2) Every slide layout's Date placeholder
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes $newDate
}
